# Automatic update of files.
#
# The underlying data for three observation records (rows 3, 4 and 6) was
# re-synced from source: row 3's record became row 4, row 4's record became
# row 6, and row 6's record became row 3 (a 3-way rotation of the record
# identity/location/observer fields), while the shared taxon/date columns
# stayed put. Apply that as direct cell writes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 gets what used to be row 4's identifying data ---
$ws.Range("A3").Value = 111454321
$ws.Range("P3").Value = "Albinvägen3, Svartnäset, Hackås, Jmt"
$ws.Range("Q3").Value = 492408.8681431987
$ws.Range("R3").Value = 6948272.081593725

# --- Row 4 gets what used to be row 6's identifying data ---
$ws.Range("A4").Value = 111454959
$ws.Range("P4").Value = "Siljebodarna, Jmt"
$ws.Range("Q4").Value = 492425.096130528
$ws.Range("R4").Value = 6948324.435442663
$ws.Range("S4").Value = 15
$ws.Range("AW4").Value = "Monica Magnesved"
$ws.Range("AX4").Value = "Monica Magnesved"
# Row 4 no longer carries these (now-blank) optional fields.
$ws.Range("J4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("AF4").ClearContents()

# --- Row 6 gets what used to be row 3's identifying data ---
$ws.Range("A6").Value = 111454300
$ws.Range("P6").Value = "Albinvägen2, Svartnäset, Hackås, Jmt"
$ws.Range("Q6").Value = 492448.9318965223
$ws.Range("R6").Value = 6948282.559996245
$ws.Range("S6").Value = 10
$ws.Range("AW6").Value = "Jan Magnesved"
$ws.Range("AX6").Value = "Jan Magnesved, Anders Wännström "
# Row 6 now gains these (blank, present) optional fields.
$ws.Range("J6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("AF6").Value = ""
